$d = $word.ActiveDocument

$replacements = @(
    @{old = "525×9="; new = "489×8="},
    @{old = "613×8="; new = "430×7="},
    @{old = "340×2="; new = "611×2="},
    @{old = "973×6="; new = "488×6="},
    @{old = "445×8="; new = "477×9="},
    @{old = "879×5="; new = "816×5="},
    @{old = "365×8="; new = "467×3="},
    @{old = "663×2="; new = "130×2="},
    @{old = "530×8="; new = "234×9="},
    @{old = "228×7="; new = "325×3="},
    @{old = "360×3="; new = "701×9="},
    @{old = "332×9="; new = "675×6="},
    @{old = "239×6="; new = "198×4="},
    @{old = "669×8="; new = "942×3="},
    @{old = "535×4="; new = "811×8="},
    @{old = "197×2="; new = "335×8="},
    @{old = "959×3="; new = "362×4="},
    @{old = "953×4="; new = "259×2="},
    @{old = "238×8="; new = "665×6="},
    @{old = "839×2="; new = "616×7="},
    @{old = "106×4="; new = "341×2="},
    @{old = "739×8="; new = "723×9="},
    @{old = "112×9="; new = "693×6="},
    @{old = "768×8="; new = "655×9="},
    @{old = "128×2="; new = "140×8="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
